# Commit code ngay 03/03/2020
# Update Schedule sheet: task "Code" sub-rows get fleshed out with real
# start/finish (planned + actual) dates, progress percentage, and task text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Row 9 ("Code" task, first line): fill planned/actual start & finish dates,
# bump progress to 90%.
$ws.Range("C9").Value = "Tạo Database, Khởi tạo, Chọn trong Combobox, Thực hiện tìm kiếm"
$ws.Range("D9").Value = 43892
$ws.Range("E9").Value = 43893
$ws.Range("F9").Value = 43892
$ws.Range("G9").Value = 43893
$ws.Range("H9").Value = 0.9

# Row 10 ("Code" task, second line): new sub task with its own dates.
$ws.Range("C10").Value = "Thêm, Sửa, Xóa, Lưu"
$ws.Range("D10").Value = 43894
$ws.Range("E10").Value = 43896

# Reflect the author's last selection in the sheet.
$ws.Range("H9").Select()
